$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1147.8334
$ws.Range("I2").Value = 1321.75
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1321.75
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -1208.75
$ws.Range("N2").Value = -1026
$ws.Range("H4").Value = 330.16666
$ws.Range("I4").Value = 196.2
$ws.Range("K4").Value = 196.2
$ws.Range("M4").Value = -82.19999999999999
$ws.Range("H9").Value = 1192.3334
$ws.Range("I9").Value = 1638.75
$ws.Range("K9").Value = 1638.75
$ws.Range("M9").Value = -1469.75
$ws.Range("H17").Value = 2158.875
$ws.Range("J17").Value = 2158.875
$ws.Range("L17").Value = 6476.625
$ws.Range("N17").Value = -6812.625
$ws.Range("H28").Value = 1351.2858
$ws.Range("I28").Value = 1077.3334
$ws.Range("K28").Value = 1077.3334
$ws.Range("M28").Value = -592.3334
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H53").Value = 667.25
$ws.Range("I53").Value = 623.1667
$ws.Range("J53").Value = 799.5
$ws.Range("K53").Value = 623.1667
$ws.Range("L53").Value = 799.5
$ws.Range("M53").Value = 13.83330000000001
$ws.Range("N53").Value = -2073.5
$ws.Range("H92").Value = 377.3846
$ws.Range("I92").Value = 325
$ws.Range("J92").Value = 665.5
$ws.Range("K92").Value = 325
$ws.Range("L92").Value = 665.5
$ws.Range("M92").Value = 923
$ws.Range("N92").Value = -3161.5
$ws.Range("H98").Value = 1322.6316
$ws.Range("I98").Value = 947.25
$ws.Range("K98").Value = 947.25
$ws.Range("M98").Value = 550.75
$ws.Range("H121").Value = 1841
$ws.Range("J121").Value = 1932.8948
$ws.Range("L121").Value = 5798.6844
$ws.Range("N121").Value = -9292.6844
$ws.Range("H122").Value = 1322.6316
$ws.Range("I122").Value = 947.25
$ws.Range("K122").Value = 2841.75
$ws.Range("M122").Value = -391.75
$ws.Range("H137").Value = 4762.364
$ws.Range("I137").Value = 3316
$ws.Range("K137").Value = 9948
$ws.Range("M137").Value = -7398

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1461.3334
$ws.Range("I61").Value = 1461.3334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1461.3334
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1249.3334
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2773.5
$ws.Range("I74").Value = 2448.3928
$ws.Range("K74").Value = 2448.3928
$ws.Range("M74").Value = -1574.3928
$ws.Range("H77").Value = 2773.5
$ws.Range("I77").Value = 2448.3928
$ws.Range("K77").Value = 12241.964
$ws.Range("M77").Value = -7873.964
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 3588.2778
$ws.Range("I132").Value = 2599.3635
$ws.Range("J132").Value = 5142.2856
$ws.Range("K132").Value = 7798.0905
$ws.Range("L132").Value = 15426.8568
$ws.Range("M132").Value = -5268.0905
$ws.Range("N132").Value = -20486.8568
$ws.Range("H136").Value = 1461.3334
$ws.Range("I136").Value = 1461.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4384.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1834.0002
$ws.Range("N136").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 927.46155
$ws.Range("J80").Value = 870
$ws.Range("L80").Value = 870
$ws.Range("N80").Value = -2866
$ws.Range("H83").Value = 927.46155
$ws.Range("J83").Value = 870
$ws.Range("L83").Value = 4350
$ws.Range("N83").Value = -14334
$ws.Range("H105").Value = 3740.2144
$ws.Range("I105").Value = 3651
$ws.Range("K105").Value = 3651
$ws.Range("M105").Value = -1904

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 155
$ws.Range("I7").Value = 115.42857
$ws.Range("K7").Value = 115.42857
$ws.Range("M7").Value = -2.428569999999993
$ws.Range("H31").Value = 4709.7144
$ws.Range("I31").Value = 4992.25
$ws.Range("J31").Value = 4333
$ws.Range("K31").Value = 4992.25
$ws.Range("L31").Value = 4333
$ws.Range("M31").Value = -4697.25
$ws.Range("N31").Value = -4923
$ws.Range("H34").Value = 4709.7144
$ws.Range("I34").Value = 4992.25
$ws.Range("J34").Value = 4333
$ws.Range("K34").Value = 4992.25
$ws.Range("L34").Value = 4333
$ws.Range("M34").Value = -4790.25
$ws.Range("N34").Value = -4737
$ws.Range("H58").Value = 2562.5715
$ws.Range("J58").Value = 3013.25
$ws.Range("L58").Value = 3013.25
$ws.Range("N58").Value = -3419.25
$ws.Range("H107").Value = 1607.3334
$ws.Range("I107").Value = 698.75
$ws.Range("K107").Value = 698.75
$ws.Range("M107").Value = 1221.25
$ws.Range("H136").Value = 2562.5715
$ws.Range("J136").Value = 3013.25
$ws.Range("L136").Value = 9039.75
$ws.Range("N136").Value = -14139.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 86781050
$ws.Range("I4").Value = 25925984
$ws.Range("J4").Value = 121555370
$ws.Range("K4").Value = 77777952
$ws.Range("L4").Value = 364666110
$ws.Range("M4").Value = -77777840
$ws.Range("N4").Value = -364666334
$ws.Range("H12").Value = 210.58824
$ws.Range("J12").Value = 178.66667
$ws.Range("L12").Value = 536.00001
$ws.Range("N12").Value = -882.00001
$ws.Range("H37").Value = 199999.67
$ws.Range("J37").Value = 199999.67
$ws.Range("L37").Value = 599999.01
$ws.Range("N37").Value = -600223.01
$ws.Range("H113").Value = 1012.5714
$ws.Range("J113").Value = 996.3333
$ws.Range("L113").Value = 2988.9999
$ws.Range("N113").Value = -7328.9999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14542.333
$ws.Range("I70").Value = 14664.667
$ws.Range("J70").Value = 14297.667
$ws.Range("K70").Value = 14664.667
$ws.Range("L70").Value = 14297.667
$ws.Range("M70").Value = -14394.667
$ws.Range("N70").Value = -14837.667
$ws.Range("H73").Value = 14542.333
$ws.Range("I73").Value = 14664.667
$ws.Range("J73").Value = 14297.667
$ws.Range("K73").Value = 14664.667
$ws.Range("L73").Value = 14297.667
$ws.Range("M73").Value = -13728.667
$ws.Range("N73").Value = -16169.667
$ws.Range("H102").Value = 2638.2222
$ws.Range("I102").Value = 955.1667
$ws.Range("J102").Value = 6004.3335
$ws.Range("K102").Value = 955.1667
$ws.Range("L102").Value = 6004.3335
$ws.Range("M102").Value = 666.8333
$ws.Range("N102").Value = -9248.333500000001
$ws.Range("H132").Value = 1164.5555
$ws.Range("I132").Value = 671.56525
$ws.Range("J132").Value = 3999.25
$ws.Range("K132").Value = 2014.69575
$ws.Range("L132").Value = 11997.75
$ws.Range("M132").Value = 515.3042500000001
$ws.Range("N132").Value = -17057.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3071
$ws.Range("I132").Value = 2352.6428
$ws.Range("K132").Value = 7057.928400000001
$ws.Range("M132").Value = -4527.928400000001
